# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose dialog-act annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4; DAMSL = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 6; DAMSL = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 12; DAMSL = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 13; DAMSL = "ba"; DialogAct = "Appreciation" }
    @{ Row = 20; DAMSL = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 22; DAMSL = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 23; DAMSL = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 24; DAMSL = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 25; DAMSL = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 28; DAMSL = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 30; DAMSL = "qy"; DialogAct = "Yes-No-Question" }
    @{ Row = 32; DAMSL = "qy"; DialogAct = "Yes-No-Question" }
    @{ Row = 35; DAMSL = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 36; DAMSL = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 44; DAMSL = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 52; DAMSL = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 57; DAMSL = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 58; DAMSL = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 68; DAMSL = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 77; DAMSL = "ba"; DialogAct = "Appreciation" }
    @{ Row = 81; DAMSL = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 85; DAMSL = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 93; DAMSL = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 99; DAMSL = "sd"; DialogAct = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSL
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
